$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# (e.g. "19.00", "0.4480") keep their exact original text instead of
# being reinterpreted as numbers and losing trailing zeros.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "30.169.29"
$ws.Range("E2").Value = "  -0.74%  "

# Row 3
$ws.Range("D3").Value = "1.907.89"
$ws.Range("E3").Value = "  -1.57%  "

# Row 4
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
$ws.Range("D5").Value = "0.7586"
$ws.Range("E5").Value = "  -1.59%  "

# Row 6
$ws.Range("D6").Value = "243.61"
$ws.Range("E6").Value = "  -2.08%  "

# Row 7
$ws.Range("D7").Value = "0.9982"
$ws.Range("E7").Value = "  -0.23%  "

# Row 8
$ws.Range("D8").Value = "0.3147"
$ws.Range("E8").Value = "  -1.84%  "

# Row 9
$ws.Range("D9").Value = "27.08"
$ws.Range("E9").Value = "  -3.24%  "

# Row 10
$ws.Range("D10").Value = "0.06967"
$ws.Range("E10").Value = "  -2.16%  "

# Row 11
$ws.Range("D11").Value = "0.7810"
$ws.Range("E11").Value = "  -0.64%  "

# Row 12
$ws.Range("D12").Value = "0.07978"
$ws.Range("E12").Value = "  -0.42%  "

# Row 13
$ws.Range("D13").Value = "1.914.02"
$ws.Range("E13").Value = "  -1.23%  "

# Row 14
$ws.Range("D14").Value = "5.278"
$ws.Range("E14").Value = "  -1.96%  "

# Row 15
$ws.Range("D15").Value = "91.78"
$ws.Range("E15").Value = "  -3.40%  "

# Row 16
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "14.30"
$ws.Range("E16").Value = "  -2.03%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "30.150.69"
$ws.Range("E17").Value = "  -0.81%  "

# Row 18
$ws.Range("D18").Value = "5.855"
$ws.Range("E18").Value = "  +0.71%  "

# Row 19
$ws.Range("D19").Value = "244.19"
$ws.Range("E19").Value = "  -5.24%  "

# Row 20
$ws.Range("D20").Value = "0.000007849"
$ws.Range("E20").Value = "  -2.25%  "

# Row 21
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "0.9986"
$ws.Range("E21").Value = "  -0.17%  "

# Row 22
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.154.83"
$ws.Range("E22").Value = "  -1.43%  "

# Row 23
$ws.Range("D23").Value = "0.9964"
$ws.Range("E23").Value = "  -0.45%  "

# Row 24
$ws.Range("D24").Value = "6.680"
$ws.Range("E24").Value = "  -1.58%  "

# Row 25
$ws.Range("D25").Value = "9.442"
$ws.Range("E25").Value = "  -1.74%  "

# Row 26
$ws.Range("D26").Value = "165.77"
$ws.Range("E26").Value = "  +0.56%  "

# Row 27
$ws.Range("D27").Value = "19.00"
$ws.Range("E27").Value = "  -0.95%  "

# Row 28
$ws.Range("D28").Value = "0.1289"
$ws.Range("E28").Value = "  -3.64%  "

# Row 29
$ws.Range("D29").Value = "2.114"
$ws.Range("E29").Value = "  -7.97%  "

# Row 30
$ws.Range("D30").Value = "1.346"
$ws.Range("E30").Value = "  -1.68%  "

# Row 31
$ws.Range("D31").Value = "1.544"
$ws.Range("E31").Value = "  +0.79%  "

# Row 32
$ws.Range("D32").Value = "4.333"
$ws.Range("E32").Value = "  -2.34%  "

# Row 33
$ws.Range("D33").Value = "4.078"
$ws.Range("E33").Value = "  -2.11%  "

# Row 34
$ws.Range("D34").Value = "0.05183"
$ws.Range("E34").Value = "  -0.43%  "

# Row 35
$ws.Range("D35").Value = "1.298"
$ws.Range("E35").Value = "  +1.29%  "

# Row 36
$ws.Range("D36").Value = "0.7465"
$ws.Range("E36").Value = "  -0.64%  "

# Row 37
$ws.Range("D37").Value = "2.753"
$ws.Range("E37").Value = "  -0.87%  "

# Row 38
$ws.Range("E38").Value = "  -1.55%  "

# Row 39
$ws.Range("D39").Value = "2.793"
$ws.Range("E39").Value = "  -0.46%  "

# Row 40
$ws.Range("D40").Value = "6.381"
$ws.Range("E40").Value = "  -1.27%  "

# Row 41
$ws.Range("D41").Value = "75.29"
$ws.Range("E41").Value = "  -3.60%  "

# Row 42
$ws.Range("D42").Value = "0.4480"
$ws.Range("E42").Value = "  -0.82%  "

# Row 43
$ws.Range("D43").Value = "1.942"
$ws.Range("E43").Value = "  -2.05%  "

# Row 44
$ws.Range("D44").Value = "0.9967"
$ws.Range("E44").Value = "  -0.47%  "

# Row 45
$ws.Range("D45").Value = "0.8371"
$ws.Range("E45").Value = "  -0.15%  "

# Row 46
$ws.Range("D46").Value = "7.639"
$ws.Range("E46").Value = "  +1.26%  "

# Row 47
$ws.Range("D47").Value = "101.23"
$ws.Range("E47").Value = "  -0.42%  "

# Row 48
$ws.Range("D48").Value = "9.830"
$ws.Range("E48").Value = "  +0.07%  "

# Row 49
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "37.43"
$ws.Range("E49").Value = "  -0.36%  "

# Row 50
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "946.34"
$ws.Range("E50").Value = "  -3.08%  "

# Row 51
$ws.Range("D51").Value = "0.1197"
$ws.Range("E51").Value = "  +1.13%  "
